$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B7 (closeStatement) to 25% complete and highlight the row label
$ws.Range("B7").Value = 0.25
$ws.Range("A7").Interior.Color = 65535

# Update B35 (readStatement) to 33% complete and highlight the row label
$ws.Range("B35").Value = 0.33
$ws.Range("A35").Interior.Color = 65535

# Move the active selection to A7
$ws.Range("A7").Select()
